$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sampdoria 2020 match log: the last six fixtures (rows 10-15) were missing
# their xG_home / xG_away / goals_home / goals_away stats. Fill them in.
$data = @{
    10 = @("1.63863", "1.38058", "2", "2")
    11 = @("0.986547", "2.66831", "1", "2")
    12 = @("1.51835", "0.713804", "2", "1")
    13 = @("1.49333", "0.779474", "1", "2")
    14 = @("2.48022", "1.04297", "3", "1")
    15 = @("1.6455", "2.16381", "2", "3")
}

$rows = 10, 11, 12, 13, 14, 15

# Columns D (xG_home), E (xG_away), F (goals_home), G (goals_away).
# The source file stores every value (even the integer goal counts) as text,
# so force text entry (NumberFormat "@") before assigning, then clear the
# format straight back off so the cell keeps the workbook's default style.
foreach ($col in 4, 5, 6, 7) {
    foreach ($row in $rows) {
        $cell = $ws.Cells.Item($row, $col)
        $cell.NumberFormat = "@"
        $cell.Value = $data[$row][$col - 4]
        $cell.ClearFormats()
    }
}
